$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: move the (hidden) "_GoBack" bookmark from the end of the
# "Tutorials" paragraph to the end of the "Functions" paragraph (right
# after the text, before the paragraph mark).
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$funcPara = $d.Paragraphs.Item(10)              # "...Functions"
$funcEnd = $funcPara.Range.End - 1               # just before the pilcrow

# NOTE: adding a bookmark with a *collapsed* range landing exactly on a
# paragraph-end boundary is mishandled by this host if done directly,
# so we pad with scratch text, drop the bookmark while it's safely
# mid-paragraph, then remove the scratch text again.
$scratch = $d.Range($funcEnd, $funcEnd)
$scratch.InsertAfter("ZZscratchZZ")

$bmRange = $d.Range($funcEnd, $funcEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

$scratchRange = $d.Range($funcEnd, $funcEnd + 11)
$scratchRange.Delete()

# ------------------------------------------------------------------
# Change 2: fold the "Memory Management" paragraph into the
# "References & Pointers" paragraph, joined by a new " & " run, so the
# paragraph reads "References & Pointers & Memory Management" as three
# separate runs.
# ------------------------------------------------------------------
$refsPara = $d.Paragraphs.Item(12)               # "References & Pointers"
$memPara = $d.Paragraphs.Item(13)                # "Memory Management"

# Insert a fresh empty paragraph right before "Memory Management" and
# give it the " & " text - this keeps it as its own run.
$splitPoint = $d.Range($memPara.Range.Start, $memPara.Range.Start)
$splitPoint.InsertParagraphBefore()

$joinPara = $d.Paragraphs.Item(13)               # the new empty paragraph
$joinText = $d.Range($joinPara.Range.Start, $joinPara.Range.Start)
$joinText.InsertAfter(" & ")

# Merge "References & Pointers" + " & " paragraphs by deleting the
# paragraph mark between them (structural merge keeps runs distinct).
$refsPara2 = $d.Paragraphs.Item(12)
$mark1 = $d.Range($refsPara2.Range.End - 1, $refsPara2.Range.End)
$mark1.Delete()

# Merge the combined paragraph with "Memory Management" the same way.
$combined = $d.Paragraphs.Item(12)
$mark2 = $d.Range($combined.Range.End - 1, $combined.Range.End)
$mark2.Delete()
